$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 swap their values in columns A, B, E, F, G, H
# (Note: use Value2 for reads, since Value getter is unreliable in this runtime;
#  Value is fine for writes.)
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cell5 = $ws.Range($col + "5")
    $cell6 = $ws.Range($col + "6")

    $val5 = $cell5.Value2
    $val6 = $cell6.Value2

    $cell5.Value = $val6
    $cell6.Value = $val5
}
